# Generate Report for Handback
# Updates the handback status for item 8aac8bf9-4905-4a48-827d-ea6216a98954
# on both the zh-cn and de-de localization status sheets.

$wb = $excel.ActiveWorkbook

$displayName = "8aac8bf9-4905-4a48-827d-ea6216a98954.md"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I7").Value = $displayName
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/35d074ae4637284ab5a7072434118e5561233048/e2e/8aac8bf9-4905-4a48-827d-ea6216a98954.md", "", "", $displayName)
$wsZh.Range("I7").Font.Underline = 2
$wsZh.Range("I7").Font.Color = 15570276

$wsZh.Range("J7").Value = "8aac8bf9-4905-4a48-827d-ea6216a98954.b41acfed4a542616130b07ce9aa77b73966a0582.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-30 22:59:58"
$wsZh.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/62c2ea27f24e80887b10fcb70b07adf3eacc3e6b/e2e/8aac8bf9-4905-4a48-827d-ea6216a98954.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/35d074ae4637284ab5a7072434118e5561233048/e2e/8aac8bf9-4905-4a48-827d-ea6216a98954.md."

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I7").Value = $displayName
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/35d074ae4637284ab5a7072434118e5561233048/e2e/8aac8bf9-4905-4a48-827d-ea6216a98954.md", "", "", $displayName)
$wsDe.Range("I7").Font.Underline = 2
$wsDe.Range("I7").Font.Color = 15570276

$wsDe.Range("J7").Value = "8aac8bf9-4905-4a48-827d-ea6216a98954.b41acfed4a542616130b07ce9aa77b73966a0582.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-30 23:00:27"
$wsDe.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/62c2ea27f24e80887b10fcb70b07adf3eacc3e6b/e2e/8aac8bf9-4905-4a48-827d-ea6216a98954.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/35d074ae4637284ab5a7072434118e5561233048/e2e/8aac8bf9-4905-4a48-827d-ea6216a98954.md."
